# Applies the "evaluation graph is completed" edit to MainDataset:
#  - renames the 6 header labels to short machine-friendly names
#  - corrects the data in rows 1181-1200 and appends new rows 1201-1220
#  - resizes columns A:F and restores the saved selection/scroll state
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header labels (row 1) ---
$ws.Range("A1").Value = "AffectedAge"
$ws.Range("B1").Value = "MarriageAgeMother"
$ws.Range("C1").Value = "AgeOfDelivery"
$ws.Range("D1").Value = "DiseaseDuringPregnancy"
$ws.Range("E1").Value = "Disease"
$ws.Range("F1").Value = "DiseaseName"

# --- Data rows 1181:1220 (existing rows corrected, 1201:1220 newly added) ---
$rows = @(
  @(0,15,18,1,1,2),
  @(2,18,15,0,1,0),
  @(0,23,16,1,1,0),
  @(1,30,22,0,0,3),
  @(0,26,26,0,0,3),
  @(2,24,26,0,0,3),
  @(1,16,20,1,1,0),
  @(0,16,16,1,1,2),
  @(2,25,25,0,1,0),
  @(2,25,23,0,1,1),
  @(0,15,18,1,1,1),
  @(2,25,19,0,1,1),
  @(0,25,24,1,0,3),
  @(0,15,18,1,1,1),
  @(2,26,26,0,0,3),
  @(1,20,24,0,1,2),
  @(2,19,17,0,1,2),
  @(1,16,15,1,1,0),
  @(1,21,19,1,1,0),
  @(0,15,17,1,1,0),
  @(1,20,24,0,1,2),
  @(2,19,18,0,1,2),
  @(1,17,23,0,1,1),
  @(1,21,19,1,1,0),
  @(0,25,24,1,0,3),
  @(0,14,16,1,1,0),
  @(1,26,27,0,0,3),
  @(1,17,24,0,1,2),
  @(1,17,24,0,1,2),
  @(0,14,15,1,1,0),
  @(2,28,28,0,0,3),
  @(1,20,24,0,1,2),
  @(0,14,16,1,1,0),
  @(0,14,15,1,1,0),
  @(0,14,16,1,1,0),
  @(0,15,18,1,1,1),
  @(1,30,22,0,0,3),
  @(1,16,19,1,1,0),
  @(1,17,24,0,1,2),
  @(0,14,16,1,1,0)
)
$data = New-Object 'object[,]' $rows.Count,6
for ($i = 0; $i -lt $rows.Count; $i++) {
  for ($j = 0; $j -lt 6; $j++) {
    $data[$i,$j] = $rows[$i][$j]
  }
}
$ws.Range("A1181:F1220").Value = $data

# --- Column widths (A:F), values chosen so the COM pixel-grid rounding
#     lands as close as possible to the widths recorded in the workbook) ---
$ws.Columns.Item(1).ColumnWidth = 26.15
$ws.Columns.Item(2).ColumnWidth = 24.35
$ws.Columns.Item(3).ColumnWidth = 17.65
$ws.Columns.Item(4).ColumnWidth = 23.35
$ws.Columns.Item(5).ColumnWidth = 15.5
$ws.Columns.Item(6).ColumnWidth = 17.65

# --- Restore view state: scrolled position + active selection ---
$win = $excel.Windows.Item(1)
$win.ScrollRow = 1165
$win.ScrollColumn = 1
$ws.Application.GoTo($ws.Range("A1181"), $false)
$ws.Range("A1181:F1220").Select() | Out-Null
